# Add a new "2021" data column (R) to the table, mirroring the layout of
# the existing last column (Q / 2020) and then updating the figures that
# are specific to the new year.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy column Q's (2020) cells for the data block (header year + the four
# data rows) into column R so the new column inherits the same number
# formats/borders/fonts used throughout the table.
$ws.Range("Q3:Q8").Copy($ws.Range("R3:R8"))

# Year label for the new column.
$ws.Range("R3").Value = 2021

# Raw 2021 figures.
$ws.Range("R6").Value = 312
$ws.Range("R7").Value = 1910
$ws.Range("R8").Value = 4409166

# Derived "per 100 000 adults" ratios, same formula pattern as every other
# year column.
$ws.Range("R4").Formula = "=R6/R8*100000"
$ws.Range("R5").Formula = "=R7/R8*100000"

# Match the saved view state: selection moved to the new column.
$ws.Range("R15").Select()
